$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "41.296.94"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.86%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.181.35"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("E4").Value = "  -0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "238.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.610"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.44%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "70.16"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.55%  "

$ws.Range("E8").Value = "  +0.04%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.580"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.50%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.17"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -6.83%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0926"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.39"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.42%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.76"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.93%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.504.16"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.00"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.57%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.800"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.82%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.170.50"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.56%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "41.148.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000101"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -6.71%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "70.57"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.99%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.95"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "

$ws.Range("E23").Value = "  -4.51%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "226.20"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.58%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.46%  "

$ws.Range("E26").Value = "  +0.11%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.71%  "

$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("E29").Value = "  -2.59%  "

$ws.Range("E30").Value = "  +0.25%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "167.95"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "19.95"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.17%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "31.12"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +7.61%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0770"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.91%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.14"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -8.64%  "

$ws.Range("E36").Value = "  -3.27%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -7.20%  "

$ws.Range("E38").Value = "  -3.62%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0286"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.94%  "

$ws.Range("E40").Value = "  -1.44%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -10.10%  "

$ws.Range("E42").Value = "  -3.17%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "59.76"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -8.17%  "

$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("E45").Value = "  -2.27%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.30"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.08%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "98.48"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.75%  "

$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("E49").Value = "  -2.59%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.21"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -8.07%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.63"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.53%  "
